# Weekly fruit/vegetable price update - Apio (Vega Modelo de Temuco)
# Two new daily records are inserted ahead of the existing tail rows
# (which shift down from 165-167 to 167-169 unchanged).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 165 (pushes old 165..167 down to 167..169)
$ws.Rows.Item(165).Insert()
$ws.Rows.Item(165).Insert()

# New row 165
$ws.Cells.Item(165, 1).Value = 10
$ws.Cells.Item(165, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(165, 3).Value = "La Araucanía"
$ws.Cells.Item(165, 4).Value = 44448
$ws.Cells.Item(165, 5).Value = 9
$ws.Cells.Item(165, 6).Value = 100112017
$ws.Cells.Item(165, 7).Value = "Apio"
$ws.Cells.Item(165, 8).Value = "Americana (o)"
$ws.Cells.Item(165, 9).Value = "Primera"
$ws.Cells.Item(165, 10).Value = 85
$ws.Cells.Item(165, 11).Value = 10000
$ws.Cells.Item(165, 12).Value = 10000
$ws.Cells.Item(165, 13).Value = 10000
$ws.Cells.Item(165, 14).Value = "$/docena de matas"
$ws.Cells.Item(165, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(165, 16).Value = 1667
$ws.Cells.Item(165, 17).Value = 6
$ws.Cells.Item(165, 18).Value = "Hortaliza"

# New row 166
$ws.Cells.Item(166, 1).Value = 10
$ws.Cells.Item(166, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(166, 3).Value = "La Araucanía"
$ws.Cells.Item(166, 4).Value = 44448
$ws.Cells.Item(166, 5).Value = 9
$ws.Cells.Item(166, 6).Value = 100112017
$ws.Cells.Item(166, 7).Value = "Apio"
$ws.Cells.Item(166, 8).Value = "Americana (o)"
$ws.Cells.Item(166, 9).Value = "Primera"
$ws.Cells.Item(166, 10).Value = 110
$ws.Cells.Item(166, 11).Value = 10000
$ws.Cells.Item(166, 12).Value = 10000
$ws.Cells.Item(166, 13).Value = 10000
$ws.Cells.Item(166, 14).Value = "$/docena de matas"
$ws.Cells.Item(166, 15).Value = "Región Metropolitana"
$ws.Cells.Item(166, 16).Value = 1667
$ws.Cells.Item(166, 17).Value = 6
$ws.Cells.Item(166, 18).Value = "Hortaliza"
